# Weekly update: insert a new data row at the top of the price table
# (row 2), pushing all existing records down by one row, and fill the
# new row with this week's Damasco (apricot) price report for the
# "Vega Modelo de Temuco" market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 2, shifting rows 2:83 down to 3:84.
# CopyOrigin = 0 (xlFormatFromLeftOrAbove) is Excel's default for this
# operation, but since the header row above is bold/bordered we strip
# the inherited formatting right after so the new data row matches the
# plain look of the rest of the table.
$ws.Rows.Item(2).Insert(-4121, 0)
$ws.Rows.Item(2).ClearFormats()

# Column D holds dates stored as plain numbers formatted as a date/time;
# re-apply that same custom number format to the new row's date cell.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44922
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Dina"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 35
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1333
$ws.Range("T2").Value = 15
